# Apply updated forecast-error values to rows 2-5, and append new rows 6-11
# (Q4..Q9) with the same structure/formatting as the existing Q0-Q3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Update existing rows 2-5 (B:G) with revised values ---
$existing = @(
    @(2, 0.06224550739541238, 0.4815100860261196, 0.5188979745415299, 0.7203457326461578, 0.7255814477924799, 46),
    @(3, 0.009644723670651665, 0.4809544651665302, 0.5144163341041081, 0.7172282301360622, 0.725267178095959, 45),
    @(4, 0.04605011742652749, 0.4545971144787497, 0.4787573890230339, 0.691922964659386, 0.6983704975779645, 44),
    @(5, 0.02000496247511335, 0.4699992430775199, 0.4788062426551278, 0.6919582665559592, 0.6998547450934846, 43)
)

foreach ($row in $existing) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# --- Append new rows 6-11 for Q4..Q9 ---
$newRows = @(
    @(6, "Q4", 0.05983179997382334, 0.4610733700736253, 0.4940437835170411, 0.702882482010358, 0.7088204810331585, 42),
    @(7, "Q5", 0.02666374005585367, 0.4762874528960522, 0.4947189357917796, 0.7033625919764142, 0.7115884908266601, 41),
    @(8, "Q6", 0.05161238290339588, 0.4630942636991037, 0.5133492645108798, 0.7164839597024345, 0.7237264155437525, 40),
    @(9, "Q7", 0.03337024079663585, 0.4821917770522179, 0.5171250100704028, 0.7191140452462341, 0.7277298162397244, 39),
    @(10, "Q8", 0.06108869437512898, 0.4568004358409793, 0.5137306737460197, 0.7167500776044741, 0.7237282623592945, 38),
    @(11, "Q9", 0.03255594249262348, 0.4650544728812123, 0.502363057293461, 0.708775745418437, 0.7177940205234036, 37)
)

# Use the formatting of the existing Q3 label cell (A5) as the template for
# the newly appended label cells (bold font, border, centered/top aligned).
$styleTemplate = $ws.Cells.Item(5, 1)

foreach ($row in $newRows) {
    $r = $row[0]

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[1]
    $styleTemplate.Copy()
    $cellA.PasteSpecial($xlPasteFormats)

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
